$changes = @(
    @{ Cell="D2"; Value="28.133.05" }
    @{ Cell="E2"; Value="  +2.17%  " }
    @{ Cell="D3"; Value="1.911.23" }
    @{ Cell="E3"; Value="  +2.09%  " }
    @{ Cell="E4"; Value="  -0.94%  " }
    @{ Cell="D5"; Value="316.76" }
    @{ Cell="E5"; Value="  +1.26%  " }
    @{ Cell="E6"; Value="  -0.97%  " }
    @{ Cell="D7"; Value="0.4821" }
    @{ Cell="E7"; Value="  +0.79%  " }
    @{ Cell="E8"; Value="  +1.10%  " }
    @{ Cell="D9"; Value="0.07369" }
    @{ Cell="E9"; Value="  +0.04%  " }
    @{ Cell="D10"; Value="0.9339" }
    @{ Cell="E10"; Value="  -0.47%  " }
    @{ Cell="D11"; Value="20.85" }
    @{ Cell="E11"; Value="  +0.49%  " }
    @{ Cell="B12"; Value="TRON" }
    @{ Cell="C12"; Value="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx" }
    @{ Cell="D12"; Value="0.07805" }
    @{ Cell="E12"; Value="  -0.62%  " }
    @{ Cell="B13"; Value="WrappedEther" }
    @{ Cell="C13"; Value="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" }
    @{ Cell="D13"; Value="1.884.46" }
    @{ Cell="E13"; Value="  +0.37%  " }
    @{ Cell="D14"; Value="5.519" }
    @{ Cell="D15"; Value="6.655" }
    @{ Cell="E15"; Value="  +0.91%  " }
    @{ Cell="D16"; Value="91.70" }
    @{ Cell="E16"; Value="  +0.84%  " }
    @{ Cell="D17"; Value="1.005" }
    @{ Cell="E17"; Value="  -1.05%  " }
    @{ Cell="D18"; Value="0.000008829" }
    @{ Cell="E18"; Value="  -1.15%  " }
    @{ Cell="E19"; Value="  -0.86%  " }
    @{ Cell="D20"; Value="28.170.06" }
    @{ Cell="E20"; Value="  +2.18%  " }
    @{ Cell="D21"; Value="14.86" }
    @{ Cell="E21"; Value="  -0.74%  " }
    @{ Cell="D22"; Value="5.172" }
    @{ Cell="E22"; Value="  +0.58%  " }
    @{ Cell="D23"; Value="2.172.43" }
    @{ Cell="E23"; Value="  +2.43%  " }
    @{ Cell="D24"; Value="10.92" }
    @{ Cell="E24"; Value="  +1.61%  " }
    @{ Cell="D25"; Value="156.05" }
    @{ Cell="E25"; Value="  +1.13%  " }
    @{ Cell="D26"; Value="1.921" }
    @{ Cell="E26"; Value="  -2.17%  " }
    @{ Cell="D27"; Value="18.55" }
    @{ Cell="E27"; Value="  -0.04%  " }
    @{ Cell="D28"; Value="2.115" }
    @{ Cell="E28"; Value="  +4.35%  " }
    @{ Cell="D29"; Value="116.52" }
    @{ Cell="E29"; Value="  +0.44%  " }
    @{ Cell="D30"; Value="4.971" }
    @{ Cell="E30"; Value="  -0.67%  " }
    @{ Cell="D31"; Value="0.08957" }
    @{ Cell="E31"; Value="  +0.22%  " }
    @{ Cell="D32"; Value="3.336" }
    @{ Cell="E32"; Value="  +0.11%  " }
    @{ Cell="D33"; Value="1.258" }
    @{ Cell="E33"; Value="  +3.23%  " }
    @{ Cell="D34"; Value="0.7747" }
    @{ Cell="E34"; Value="  +2.49%  " }
    @{ Cell="D35"; Value="4.694" }
    @{ Cell="E35"; Value="  +1.73%  " }
    @{ Cell="D36"; Value="2.648" }
    @{ Cell="E36"; Value="  -1.83%  " }
    @{ Cell="D37"; Value="0.02061" }
    @{ Cell="E37"; Value="  +0.20%  " }
    @{ Cell="D38"; Value="1.108" }
    @{ Cell="E38"; Value="  -0.94%  " }
    @{ Cell="D39"; Value="0.05326" }
    @{ Cell="E39"; Value="  +0.66%  " }
    @{ Cell="D40"; Value="0.5501" }
    @{ Cell="E40"; Value="  +2.75%  " }
    @{ Cell="D41"; Value="2.987" }
    @{ Cell="E41"; Value="  -0.51%  " }
    @{ Cell="D42"; Value="7.029" }
    @{ Cell="E42"; Value="  -0.88%  " }
    @{ Cell="D43"; Value="0.1531" }
    @{ Cell="E43"; Value="  +0.25%  " }
    @{ Cell="D44"; Value="8.494" }
    @{ Cell="E44"; Value="  +0.39%  " }
    @{ Cell="D45"; Value="10.74" }
    @{ Cell="E45"; Value="  +0.36%  " }
    @{ Cell="D46"; Value="0.4843" }
    @{ Cell="E46"; Value="  +0.37%  " }
    @{ Cell="D47"; Value="108.14" }
    @{ Cell="E47"; Value="  +4.99%  " }
    @{ Cell="D48"; Value="1.004" }
    @{ Cell="E48"; Value="  -1.04%  " }
    @{ Cell="D49"; Value="1.658" }
    @{ Cell="E49"; Value="  -0.29%  " }
    @{ Cell="D50"; Value="68.58" }
    @{ Cell="E50"; Value="  +1.69%  " }
    @{ Cell="E51"; Value="  -0.18%  " }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($change in $changes) {
    $cellRef = $change.Cell
    $value = $change.Value
    $isNumeric = $value -match '^[+-]?\d+(\.\d+)?$'

    $range = $ws.Range($cellRef)
    if ($isNumeric) {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
